# working_hours.xlsx – add the missing time-tracking entry for 2014-03-29
# (a short ~1h stint that was left out), which pushes the summary block
# (sum [min] / sum [h] / sum [working weeks]) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, shifting the old row 102 (blank spacer row)
# and the three summary rows down to 103-106. Excel auto-adjusts the
# SUM(F2:F102) formula (and every other reference) to SUM(F2:F103), etc.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row with the extra entry.
$ws.Range("A102").Value = 2014
$ws.Range("B102").Value = 3
$ws.Range("C102").Value = 29
$ws.Range("D102").Value = 0.83333333333333337
$ws.Range("E102").Value = 0.875
$ws.Range("F102").Formula = "=(E102-D102)*24*60"
$ws.Range("G102").Formula = "=F102/60"

# Move the active selection to A103, mirroring where the cursor ends up
# after inserting/filling the new row.
$null = $ws.Range("A103").Select()
